# Fix the combined neighborhood population and school enrollment dashboard.
# The "elementary" sheet had a handful of stale / mis-typed attendance-area
# labels in column B; correct them so they read the same as the other
# sheets/rows that reference the same schools.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("elementary")

# Row 2 (23rd & Union-Jackson): "Thurgood Marshall" -> "T. Marshall"
$ws.Range("B2").Value = "Bailey Gatzert, Leschi, T. Marshall"

# Row 13 (Green Lake): "Greenlake" -> "Green Lake"
$ws.Range("B13").Value = "Green Lake"

# Row 18 (Mt Baker): "Thurgood Marshall" -> "T. Marshall"
$ws.Range("B18").Value = "T. Marshall, John Muir, Beacon Hill Intl, Kimball, Leschi"

# Row 23 (Roosevelt): "Greenlake" -> "Green Lake"
$ws.Range("B23").Value = "Green Lake"

# Row 30 (West Seattle Junction): remove stray space before the comma
$ws.Range("B30").Value = "Genesee Hill, Fairmount Park, Pathfinder"

# Move the active selection to B3, matching the saved view state.
$ws.Range("B3").Select()
